# Updates the "cryptos" list: refresh Price (col D) and Volume(1h) (col E)
# values for rows 2-51, matching the 2024-11-13 GitHub Actions data refresh.
#
# Note: column D values that look like plain numbers (e.g. "206.68") are
# written with a leading apostrophe so Excel keeps them as text (matching
# the workbook's original inline-string / text representation) instead of
# silently converting them to numeric cells. Values that already contain
# two dots (e.g. "3.168.88") or the "%" rows in column E are unambiguous
# text and don't need the apostrophe.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '87.649.93'
$ws.Range("E2").Value = '  -1.42%  '
$ws.Range("D3").Value = '3.168.88'
$ws.Range("E3").Value = '  -6.01%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''206.68'
$ws.Range("E5").Value = '  -7.61%  '
$ws.Range("D6").Value = '''610.74'
$ws.Range("E6").Value = '  -6.46%  '
$ws.Range("D7").Value = '''0.383'
$ws.Range("E7").Value = '  -7.47%  '
$ws.Range("D8").Value = '''0.673'
$ws.Range("E8").Value = '  +4.29%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = '3.163.92'
$ws.Range("E10").Value = '  -6.07%  '
$ws.Range("D11").Value = '''0.537'
$ws.Range("E11").Value = '  -15.21%  '
$ws.Range("D12").Value = '''0.177'
$ws.Range("E12").Value = '  +5.02%  '
$ws.Range("D13").Value = '''0.0000244'
$ws.Range("E13").Value = '  -15.29%  '
$ws.Range("D14").Value = '3.750.01'
$ws.Range("E14").Value = '  -6.07%  '
$ws.Range("D15").Value = '''5.27'
$ws.Range("E15").Value = '  -5.87%  '
$ws.Range("D16").Value = '87.363.31'
$ws.Range("E16").Value = '  -1.59%  '
$ws.Range("D17").Value = '''32.17'
$ws.Range("E17").Value = '  -13.15%  '
$ws.Range("D18").Value = '3.150.39'
$ws.Range("E18").Value = '  -6.27%  '
$ws.Range("D19").Value = '''3.02'
$ws.Range("E19").Value = '  -3.75%  '
$ws.Range("D20").Value = '''13.42'
$ws.Range("E20").Value = '  -10.51%  '
$ws.Range("D21").Value = '''417.07'
$ws.Range("E21").Value = '  -9.32%  '
$ws.Range("D22").Value = '''8.51'
$ws.Range("E22").Value = '  -12.42%  '
$ws.Range("D23").Value = '''5.11'
$ws.Range("E23").Value = '  -10.27%  '
$ws.Range("D24").Value = '''5.21'
$ws.Range("E24").Value = '  -7.52%  '
$ws.Range("D25").Value = '''11.89'
$ws.Range("E25").Value = '  -7.52%  '
$ws.Range("D26").Value = '3.329.09'
$ws.Range("E26").Value = '  -6.39%  '
$ws.Range("D27").Value = '''73.60'
$ws.Range("E27").Value = '  -8.68%  '
$ws.Range("D28").Value = '''0.0000131'
$ws.Range("E28").Value = '  -8.91%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").Value = '''0.161'
$ws.Range("E30").Value = '  -17.89%  '
$ws.Range("D31").Value = '''0.998'
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("D32").Value = '''543.55'
$ws.Range("E32").Value = '  -8.38%  '
$ws.Range("D33").Value = '''8.27'
$ws.Range("E33").Value = '  -12.45%  '
$ws.Range("E34").Value = '  -17.14%  '
$ws.Range("D35").Value = '''6.75'
$ws.Range("E35").Value = '  -6.89%  '
$ws.Range("D36").Value = '''1.85'
$ws.Range("E36").Value = '  -13.30%  '
$ws.Range("D37").Value = '''0.132'
$ws.Range("E37").Value = '  -8.84%  '
$ws.Range("D38").Value = '''21.91'
$ws.Range("E38").Value = '  -7.66%  '
$ws.Range("D39").Value = '''21.80'
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("D40").Value = '''0.999'
$ws.Range("E40").Value = '  +0.16%  '
$ws.Range("D41").Value = '''2.99'
$ws.Range("E41").Value = '  -6.10%  '
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("D43").Value = '''1.90'
$ws.Range("E43").Value = '  -12.07%  '
$ws.Range("D44").Value = '''0.371'
$ws.Range("E44").Value = '  -15.04%  '
$ws.Range("D45").Value = '''147.17'
$ws.Range("E45").Value = '  -7.14%  '
$ws.Range("D46").Value = '''173.47'
$ws.Range("E46").Value = '  -8.42%  '
$ws.Range("D47").Value = '''43.28'
$ws.Range("E47").Value = '  -7.11%  '
$ws.Range("D48").Value = '''0.127'
$ws.Range("E48").Value = '  +1.32%  '
$ws.Range("E49").Value = '  -14.50%  '
$ws.Range("D50").Value = '''3.98'
$ws.Range("E50").Value = '  -12.36%  '
$ws.Range("D51").Value = '''0.699'
$ws.Range("E51").Value = '  -11.72%  '
